$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Original state:
#   A1 = "Nombre de logs" (shared string), B1 = 0
#
# Target state:
#   A1 = "Nb log", B1 = 0
#   A2 = "Type",   B2 = "message"
#
# Apply in an order that reproduces the same shared-string table as the
# authored edit: the existing string slot (originally "Nombre de logs")
# is repurposed in place to "Type", two new strings ("Nb log", "message")
# are appended, then A1/A2 are swapped so the final cell contents match.

$ws.Range("A1").Value = "Type"
$ws.Range("A2").Value = "Nb log"
$ws.Range("B2").Value = "message"

$a1 = $ws.Range("A1").Value2
$a2 = $ws.Range("A2").Value2
$ws.Range("A1").Value = $a2
$ws.Range("A2").Value = $a1

$ws.Range("B1").Value = 0

$ws.Range("C2").Select()
